$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.063.53"
$ws.Range("E2").Value = "'  -2.01%  "
$ws.Range("D3").Value = "'2.636.20"
$ws.Range("E3").Value = "'  -3.07%  "
$ws.Range("E4").Value = "'  +0.00%  "
$ws.Range("D5").Value = "'589.76"
$ws.Range("E5").Value = "'  -3.13%  "
$ws.Range("D6").Value = "'165.77"
$ws.Range("E6").Value = "'  -1.05%  "
$ws.Range("E7").Value = "'  +0.02%  "
$ws.Range("D8").Value = "'0.542"
$ws.Range("E8").Value = "'  -2.02%  "
$ws.Range("D9").Value = "'2.637.62"
$ws.Range("E9").Value = "'  -2.94%  "
$ws.Range("E10").Value = "'  -0.83%  "
$ws.Range("D11").Value = "'0.159"
$ws.Range("E11").Value = "'  +1.50%  "
$ws.Range("E12").Value = "'  -0.81%  "
$ws.Range("E13").Value = "'  -1.11%  "
$ws.Range("D14").Value = "'27.67"
$ws.Range("E14").Value = "'  -3.19%  "
$ws.Range("D15").Value = "'3.117.38"
$ws.Range("E15").Value = "'  -3.12%  "
$ws.Range("E16").Value = "'  -3.65%  "
$ws.Range("D17").Value = "'67.219.09"
$ws.Range("E17").Value = "'  -1.71%  "
$ws.Range("D18").Value = "'2.639.26"
$ws.Range("E18").Value = "'  -3.29%  "
$ws.Range("D19").Value = "'12.03"
$ws.Range("E19").Value = "'  +1.16%  "
$ws.Range("D20").Value = "'8.13"
$ws.Range("E20").Value = "'  +6.27%  "
$ws.Range("D21").Value = "'360.90"
$ws.Range("E21").Value = "'  -3.06%  "
$ws.Range("E22").Value = "'  -3.45%  "
$ws.Range("D23").Value = "'4.75"
$ws.Range("E23").Value = "'  -4.72%  "
$ws.Range("D24").Value = "'11.07"
$ws.Range("E24").Value = "'  +9.43%  "
$ws.Range("D25").Value = "'1.98"
$ws.Range("E25").Value = "'  -5.39%  "
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("D27").Value = "'70.78"
$ws.Range("E27").Value = "'  -3.19%  "
$ws.Range("E29").Value = "'  +0.17%  "
$ws.Range("D31").Value = "'552.27"
$ws.Range("E31").Value = "'  -5.95%  "
$ws.Range("D32").Value = "'7.95"
$ws.Range("E32").Value = "'  -2.95%  "
$ws.Range("D33").Value = "'1.37"
$ws.Range("E33").Value = "'  -4.11%  "
$ws.Range("D34").Value = "'1.90"
$ws.Range("E34").Value = "'  -4.31%  "
$ws.Range("E35").Value = "'  +2.54%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  +0.10%  "
$ws.Range("E37").Value = "'  -4.78%  "
$ws.Range("D38").Value = "'157.77"
$ws.Range("E38").Value = "'  -3.07%  "
$ws.Range("D39").Value = "'19.23"
$ws.Range("E39").Value = "'  -3.33%  "
$ws.Range("D40").Value = "'0.369"
$ws.Range("E40").Value = "'  -2.48%  "
$ws.Range("D41").Value = "'5.25"
$ws.Range("E41").Value = "'  -3.27%  "
$ws.Range("E43").Value = "'  -0.53%  "
$ws.Range("E44").Value = "'  +0.01%  "
$ws.Range("E45").Value = "'  -4.99%  "
$ws.Range("E46").Value = "'  -2.00%  "
$ws.Range("D47").Value = "'0.591"
$ws.Range("E47").Value = "'  -1.17%  "
$ws.Range("D48").Value = "'0.0₆0296"
$ws.Range("E48").Value = "'  -4.52%  "
$ws.Range("D49").Value = "'152.70"
$ws.Range("D50").Value = "'3.84"
$ws.Range("E50").Value = "'  -2.27%  "
$ws.Range("D51").Value = "'1.72"
$ws.Range("E51").Value = "'  -3.14%  "
